# Applies the cryptos.xlsx data refresh described by the commit:
# "Updated cryptos list on Wed Dec  6 09:52:56 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="43.955.53"'
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Value = '  +5.25%  '
$ws.Range('D3').Formula = '="2.280.26"'
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Value = '  +3.18%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Formula = '="234.35"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('E6').Value = '  +3.91%  '
$ws.Range('D7').Formula = '="66.03"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +9.39%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +7.31%  '
$ws.Range('D10').Formula = '="0.103"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +16.12%  '
$ws.Range('D11').Formula = '="57.66"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').Formula = '="26.40"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +18.94%  '
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Formula = '="2.620.28"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +3.28%  '
$ws.Range('D15').Formula = '="15.98"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('E16').Value = '  +5.40%  '
$ws.Range('D17').Formula = '="0.835"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  +4.88%  '
$ws.Range('D18').Formula = '="2.283.09"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +3.24%  '
$ws.Range('D19').Formula = '="43.744.31"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +4.90%  '
$ws.Range('D20').Formula = '="0.0₃0989"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +9.72%  '
$ws.Range('D21').Formula = '="74.14"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +2.71%  '
$ws.Range('D22').Formula = '="6.16"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').Formula = '="262.62"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  +8.26%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Formula = '="2.51"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +6.59%  '
$ws.Range('D26').Formula = '="2.32"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  -5.92%  '
$ws.Range('D27').Formula = '="10.25"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +5.84%  '
$ws.Range('D28').Formula = '="172.62"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +1.95%  '
$ws.Range('D29').Formula = '="21.09"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +6.76%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('E32').Value = '  +7.61%  '
$ws.Range('E33').Value = '  +2.79%  '
$ws.Range('D34').Formula = '="0.0689"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +6.21%  '
$ws.Range('D35').Formula = '="5.06"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('D36').Formula = '="4.76"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +2.69%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').Formula = '="6.85"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +8.34%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Formula = '="3.86"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +8.98%  '
$ws.Range('D39').Formula = '="2.38"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').Formula = '="0.0249"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +4.04%  '
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('E42').Value = '  -1.83%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').Formula = '="0.0982"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +2.94%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Formula = '="17.54"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +6.96%  '
$ws.Range('D45').Formula = '="4.48"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('D46').Formula = '="98.56"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Formula = '="1.20"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('B48').Value = 'Celestia'
$ws.Range('C48').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D48').Formula = '="10.32"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +21.16%  '
$ws.Range('D49').Formula = '="1.476.54"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('E50').Value = '  +6.70%  '
$ws.Range('D51').Formula = '="0.000206"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  -14.02%  '
$excel.CutCopyMode = $false
